$wb = $excel.ActiveWorkbook

# --- Units sheet: just select B2 (no data change) ---
$wsUnits = $wb.Worksheets.Item("Units")
$wsUnits.Select()
$wsUnits.Range("B2").Select()

# --- Demand sheet: add columns D (Solar_Plant_Kasso / unit / unit_availibility_factor)
#     and E (Power_Wholesale / node / tax_out_unit_flow) ---
$wsDemand = $wb.Worksheets.Item("Demand")

$wsDemand.Range("D1").Value = "Solar_Plant_Kasso"
$wsDemand.Range("E1").Value = "Power_Wholesale"

$wsDemand.Range("D2").Value = "unit"
$wsDemand.Range("E2").Value = "node"

$wsDemand.Range("D3").Value = "unit_availibility_factor"
$wsDemand.Range("E3").Value = "tax_out_unit_flow"

$wsDemand.Range("D4").Value = 0
$wsDemand.Range("E4").Value = 10

$wsDemand.Range("D5").Value = 0.2
$wsDemand.Range("E5").Value = 10

$wsDemand.Range("D6").Value = 0.8
$wsDemand.Range("E6").Value = 10

$wsDemand.Range("D7").Value = 0.3
$wsDemand.Range("E7").Value = 10

$wsDemand.Range("D8").Value = 0
$wsDemand.Range("E8").Value = 10

$wsDemand.Select()
$wsDemand.Range("E4").Select()
